# Apply the weekly update: two new price observations are inserted at the
# top of the "Feria Lagunitas de Puerto Montt - Cilantro" data block
# (rows 443-444), pushing the previously-existing rows (443-501) down to
# (445-503). The rest of the sheet (rows 1-442) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 443; this shifts all
# the rows below (443..501) down by two, preserving their data/styles.
$ws.Rows.Item(443).Insert()
$ws.Rows.Item(443).Insert()

# --- New row 443 -----------------------------------------------------
$ws.Cells.Item(443, 1).Value = 4
$ws.Cells.Item(443, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(443, 3).Value = "Los Lagos"
$ws.Cells.Item(443, 4).Value = 45127
$ws.Cells.Item(443, 5).Value = 10
$ws.Cells.Item(443, 6).Value = 100112040
$ws.Cells.Item(443, 7).Value = "Cilantro"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 120
$ws.Cells.Item(443, 11).Value = 12000
$ws.Cells.Item(443, 12).Value = 12000
$ws.Cells.Item(443, 13).Value = 12000
$ws.Cells.Item(443, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(443, 15).Value = "Región Metropolitana"
$ws.Cells.Item(443, 16).Value = 333
$ws.Cells.Item(443, 17).Value = 36
$ws.Cells.Item(443, 18).Value = "Hortaliza"

# --- New row 444 -----------------------------------------------------
$ws.Cells.Item(444, 1).Value = 4
$ws.Cells.Item(444, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(444, 3).Value = "Los Lagos"
$ws.Cells.Item(444, 4).Value = 45127
$ws.Cells.Item(444, 5).Value = 10
$ws.Cells.Item(444, 6).Value = 100112040
$ws.Cells.Item(444, 7).Value = "Cilantro"
$ws.Cells.Item(444, 8).Value = "Sin especificar"
$ws.Cells.Item(444, 9).Value = "Primera"
$ws.Cells.Item(444, 10).Value = 50
$ws.Cells.Item(444, 11).Value = 6000
$ws.Cells.Item(444, 12).Value = 6000
$ws.Cells.Item(444, 13).Value = 6000
$ws.Cells.Item(444, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(444, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(444, 16).Value = 3000
$ws.Cells.Item(444, 17).Value = 2
$ws.Cells.Item(444, 18).Value = "Hortaliza"
